$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (Start, End) before the existing "Time spent (hrs)" column.
# This shifts the old column G ("Time spent (hrs)") to column I.
$ws.Columns("G:H").Insert()

# Match the new columns' width to column F's width (16.54296875 raw OOXML width).
$ws.Columns("G:H").ColumnWidth = $ws.Columns("F").ColumnWidth

# New header cells for the inserted columns.
$ws.Range("G4").Value = "Start"
$ws.Range("H4").Value = "End"
$ws.Range("G4:H4").Font.Bold = $true

# Add the new time entry row.
$ws.Range("D6").Value = $ws.Range("D5").Value2
$ws.Range("E6").Value = 43365
$ws.Range("E6").NumberFormat = $ws.Range("E5").NumberFormat
$ws.Range("F6").Value = $ws.Range("F5").Value2
$ws.Range("G6").Value = 10.25
$ws.Range("H6").Value = 12.5
$ws.Range("I6").Formula = "=H6-G6"

# Update the active selection as recorded by the workbook.
[void]$ws.Range("H10").Select()
